{"js": "// Update the division-problem answers in the first (and only) table.\n// Replacements are addressed by absolute cell position (row, column) so\n// that duplicate source strings (e.g. \"25\u00f74=6, 1\" appears twice but maps\n// to two different results) are each updated independently and correctly.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// 0-indexed (row, column) -> new text. Data lives in table rows 0, 4, 8,\n// 12, 16 (the other rows are blank spacer rows), 5 columns each.\nconst updates = [\n  { row: 0, col: 0, text: \"94\u00f79=10, 4\" },\n  { row: 0, col: 1, text: \"80\u00f78=10, 0\" },\n  { row: 0, col: 2, text: \"35\u00f73=11, 2\" },\n  { row: 0, col: 3, text: \"44\u00f78=5, 4\" },\n  { row: 0, col: 4, text: \"80\u00f78=10, 0\" },\n\n  { row: 4, col: 0, text: \"14\u00f74=3, 2\" },\n  { row: 4, col: 1, text: \"74\u00f78=9, 2\" },\n  { row: 4, col: 2, text: \"88\u00f75=17, 3\" },\n  { row: 4, col: 3, text: \"18\u00f79=2, 0\" },\n  { row: 4, col: 4, text: \"39\u00f75=7, 4\" },\n\n  { row: 8, col: 0, text: \"31\u00f74=7, 3\" },\n  { row: 8, col: 1, text: \"15\u00f73=5, 0\" },\n  { row: 8, col: 2, text: \"16\u00f74=4, 0\" },\n  { row: 8, col: 3, text: \"66\u00f75=13, 1\" },\n  { row: 8, col: 4, text: \"31\u00f74=7, 3\" },\n\n  { row: 12, col: 0, text: \"47\u00f73=15, 2\" },\n  { row: 12, col: 1, text: \"30\u00f74=7, 2\" },\n  { row: 12, col: 2, text: \"60\u00f76=10, 0\" },\n  { row: 12, col: 3, text: \"65\u00f76=10, 5\" },\n  { row: 12, col: 4, text: \"73\u00f79=8, 1\" },\n\n  { row: 16, col: 0, text: \"95\u00f79=10, 5\" },\n  { row: 16, col: 1, text: \"88\u00f72=44, 0\" },\n  { row: 16, col: 2, text: \"83\u00f78=10, 3\" },\n  { row: 16, col: 3, text: \"81\u00f72=40, 1\" },\n  { row: 16, col: 4, text: \"99\u00f74=24, 3\" },\n];\n\nfor (const u of updates) {\n  const cell = table.getCell(u.row, u.col);\n  cell.value = u.text;\n}\n\nawait context.sync();\n", "ps1": "# Update the division-problem answers in the first (and only) table.\n# Replacements are addressed by absolute cell position (row, column) so\n# that duplicate source strings (e.g. \"25\u00f74=6, 1\" appears twice but maps\n# to two different results) are each updated independently and correctly.\n\n$d = $word.ActiveDocument\n$tbl = $d.Tables(1)\n\n# 1-indexed (row, column) -> new text. Data lives in table rows 1, 5, 9,\n# 13, 17 (the other rows are blank spacer rows), 5 columns each.\n$updates = @(\n    @{ Row = 1;  Col = 1; Text = \"94\u00f79=10, 4\" },\n    @{ Row = 1;  Col = 2; Text = \"80\u00f78=10, 0\" },\n    @{ Row = 1;  Col = 3; Text = \"35\u00f73=11, 2\" },\n    @{ Row = 1;  Col = 4; Text = \"44\u00f78=5, 4\" },\n    @{ Row = 1;  Col = 5; Text = \"80\u00f78=10, 0\" },\n\n    @{ Row = 5;  Col = 1; Text = \"14\u00f74=3, 2\" },\n    @{ Row = 5;  Col = 2; Text = \"74\u00f78=9, 2\" },\n    @{ Row = 5;  Col = 3; Text = \"88\u00f75=17, 3\" },\n    @{ Row = 5;  Col = 4; Text = \"18\u00f79=2, 0\" },\n    @{ Row = 5;  Col = 5; Text = \"39\u00f75=7, 4\" },\n\n    @{ Row = 9;  Col = 1; Text = \"31\u00f74=7, 3\" },\n    @{ Row = 9;  Col = 2; Text = \"15\u00f73=5, 0\" },\n    @{ Row = 9;  Col = 3; Text = \"16\u00f74=4, 0\" },\n    @{ Row = 9;  Col = 4; Text = \"66\u00f75=13, 1\" },\n    @{ Row = 9;  Col = 5; Text = \"31\u00f74=7, 3\" },\n\n    @{ Row = 13; Col = 1; Text = \"47\u00f73=15, 2\" },\n    @{ Row = 13; Col = 2; Text = \"30\u00f74=7, 2\" },\n    @{ Row = 13; Col = 3; Text = \"60\u00f76=10, 0\" },\n    @{ Row = 13; Col = 4; Text = \"65\u00f76=10, 5\" },\n    @{ Row = 13; Col = 5; Text = \"73\u00f79=8, 1\" },\n\n    @{ Row = 17; Col = 1; Text = \"95\u00f79=10, 5\" },\n    @{ Row = 17; Col = 2; Text = \"88\u00f72=44, 0\" },\n    @{ Row = 17; Col = 3; Text = \"83\u00f78=10, 3\" },\n    @{ Row = 17; Col = 4; Text = \"81\u00f72=40, 1\" },\n    @{ Row = 17; Col = 5; Text = \"99\u00f74=24, 3\" }\n)\n\nforeach ($u in $updates) {\n    $cell = $tbl.Cell($u.Row, $u.Col)\n    $rng = $cell.Range\n    # Trim the trailing end-of-cell marker so only the visible text is replaced.\n    $rng.MoveEnd(1, -1) | Out-Null\n    $rng.Text = $u.Text\n}\n"}
